# Camp Committee report bugfix: the committee-member flag (column E,
# "IsCampComittee") wasn't populated for most rows, so reports couldn't be
# generated. Replace the single placeholder data row with the real set of
# camp-committee / participant rows pulled from the source system.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ID=1) is untouched. Row 3 (the old ID=12 placeholder row) is
# removed entirely.
$ws.Range("A3:E3").Clear()

# New data rows, written at their original (non-contiguous) row numbers.
$ws.Range("A5").Value = 5
$ws.Range("B5").Value = 8
$ws.Range("C5").Value = 11
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = $true

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 14
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = $true

$ws.Range("A12").Value = 11
$ws.Range("B12").Value = 10
$ws.Range("C12").Value = 7
$ws.Range("D12").Value = 2
$ws.Range("E12").Value = $true

$ws.Range("A13").Value = 12
$ws.Range("B13").Value = 11
$ws.Range("C13").Value = 14
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = $false

$ws.Range("A14").Value = 13
$ws.Range("B14").Value = 12
$ws.Range("C14").Value = 14
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = $false

$ws.Range("A15").Value = 14
$ws.Range("B15").Value = 2
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = $false

# Match the author's final selection/cursor position.
$ws.Range("E14").Select()
